$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 7,41

# Row 25
$arr[0,0] = "WGE 340"
$arr[0,1] = "Western Interior Designers & Marine Contractors"
$arr[0,2] = "16-01-2026"
$arr[0,3] = 286962
$arr[0,4] = "Western Interior Designers & Marine Contractors"
$arr[0,5] = 34413429360
$arr[0,6] = "NEFT"
$arr[0,7] = "SBIN0003229"
$arr[0,8] = "AAAFW8862C"
$arr[0,9] = "32AAAFW8862C1Z9"
$arr[0,10] = "SAMAL KUMAR"
$arr[0,11] = "6c398193-a5d6-465d-9ff7-ff1a53d234e6"
$arr[0,12] = 42854727464
$arr[0,13] = "SBIN0060181"
$arr[0,14] = ""
$arr[0,15] = ""
$arr[0,16] = ""
$arr[0,17] = ""
$arr[0,18] = ""
$arr[0,19] = ""
$arr[0,20] = "pending"
$arr[0,21] = 15000
$arr[0,22] = ""
$arr[0,23] = "Salary Advance RPA_UNIQUE_ID : cefaf110-3962-46fe-aa83-470d9a2fda66"
$arr[0,24] = "HO"
$arr[0,25] = "URGENT"
$arr[0,26] = "hrm@westernidc.com"
$arr[0,27] = "ESTIMATION NOT MATCHED"
$arr[0,28] = 0
$arr[0,29] = 0
$arr[0,30] = 0
$arr[0,31] = ""
$arr[0,32] = ""
$arr[0,33] = ""
$arr[0,34] = ""
$arr[0,35] = ""
$arr[0,36] = ""
$arr[0,37] = ""
$arr[0,38] = ""
$arr[0,39] = ""
$arr[0,40] = ""

# Row 26
$arr[1,0] = "WGE 195"
$arr[1,1] = "Western Interior Designers & Marine Contractors"
$arr[1,2] = "16-01-2026"
$arr[1,3] = 286962
$arr[1,4] = "Western Interior Designers & Marine Contractors"
$arr[1,5] = 34413429360
$arr[1,6] = "NEFT"
$arr[1,7] = "SBIN0003229"
$arr[1,8] = "AAAFW8862C"
$arr[1,9] = "32AAAFW8862C1Z9"
$arr[1,10] = "SAYAN BATTACHARYA"
$arr[1,11] = "7a6cd1c3-f357-4acf-8b9a-de86fbde7008"
$arr[1,12] = 14810110034736
$arr[1,13] = "UCBA0001481"
$arr[1,14] = ""
$arr[1,15] = ""
$arr[1,16] = ""
$arr[1,17] = ""
$arr[1,18] = ""
$arr[1,19] = ""
$arr[1,20] = "pending"
$arr[1,21] = 1500
$arr[1,22] = ""
$arr[1,23] = "Travel Expense (Kolkatta to GOA) RPA_UNIQUE_ID : 424486a2-49d3-409c-aaf2-e9d10c42e435"
$arr[1,24] = "GOA"
$arr[1,25] = 0
$arr[1,26] = "hrm@westernidc.com"
$arr[1,27] = "ESTIMATION NOT MATCHED"
$arr[1,28] = 0
$arr[1,29] = 0
$arr[1,30] = 0
$arr[1,31] = ""
$arr[1,32] = ""
$arr[1,33] = ""
$arr[1,34] = ""
$arr[1,35] = ""
$arr[1,36] = ""
$arr[1,37] = ""
$arr[1,38] = ""
$arr[1,39] = ""
$arr[1,40] = ""

# Row 27
$arr[2,0] = "WGE 318"
$arr[2,1] = "Western Interior Designers & Marine Contractors"
$arr[2,2] = "16-01-2026"
$arr[2,3] = 286962
$arr[2,4] = "Western Interior Designers & Marine Contractors"
$arr[2,5] = 34413429360
$arr[2,6] = "NEFT"
$arr[2,7] = "SBIN0003229"
$arr[2,8] = "AAAFW8862C"
$arr[2,9] = "32AAAFW8862C1Z9"
$arr[2,10] = "KURIAKOSE GEORGE"
$arr[2,11] = "52e43895-0a0f-4ed1-8538-bd0fab14cf39"
$arr[2,12] = 38067112026
$arr[2,13] = "SBIN0071006"
$arr[2,14] = ""
$arr[2,15] = ""
$arr[2,16] = ""
$arr[2,17] = ""
$arr[2,18] = ""
$arr[2,19] = ""
$arr[2,20] = "pending"
$arr[2,21] = 1200
$arr[2,22] = ""
$arr[2,23] = "Own expense Reimbursement RPA_UNIQUE_ID : a5169d2b-47f4-4380-9d4b-7492256adb95"
$arr[2,24] = "Bellari Project"
$arr[2,25] = 0
$arr[2,26] = "hrm@westernidc.com"
$arr[2,27] = "ESTIMATION NOT MATCHED"
$arr[2,28] = 0
$arr[2,29] = 0
$arr[2,30] = 0
$arr[2,31] = ""
$arr[2,32] = ""
$arr[2,33] = ""
$arr[2,34] = ""
$arr[2,35] = ""
$arr[2,36] = ""
$arr[2,37] = ""
$arr[2,38] = ""
$arr[2,39] = ""
$arr[2,40] = ""

# Row 28
$arr[3,0] = "WGE 10"
$arr[3,1] = "Western Interior Designers & Marine Contractors"
$arr[3,2] = "16-01-2026"
$arr[3,3] = 286962
$arr[3,4] = "Western Interior Designers & Marine Contractors"
$arr[3,5] = 34413429360
$arr[3,6] = "DCR"
$arr[3,7] = "SBIN0003229"
$arr[3,8] = "AAAFW8862C"
$arr[3,9] = "32AAAFW8862C1Z9"
$arr[3,10] = "KIRAN KUMAR K"
$arr[3,11] = "bcda3e5f-7017-4035-bcb4-0d422e7b11cb"
$arr[3,12] = 30060475288
$arr[3,13] = "SBIN0009122"
$arr[3,14] = ""
$arr[3,15] = ""
$arr[3,16] = ""
$arr[3,17] = ""
$arr[3,18] = ""
$arr[3,19] = ""
$arr[3,20] = "pending"
$arr[3,21] = 2000
$arr[3,22] = ""
$arr[3,23] = "Casual Payment (Leave compenstate  09.01.2026 ) (Aneesh Mohanan, Udayan) RPA_UNIQUE_ID : 827af3af-38a1-45a8-9e54-71188839fc7c"
$arr[3,24] = "IOCL Willington"
$arr[3,25] = 0
$arr[3,26] = "hrm@westernidc.com"
$arr[3,27] = "ESTIMATION NOT MATCHED"
$arr[3,28] = 0
$arr[3,29] = 0
$arr[3,30] = 0
$arr[3,31] = ""
$arr[3,32] = ""
$arr[3,33] = ""
$arr[3,34] = ""
$arr[3,35] = ""
$arr[3,36] = ""
$arr[3,37] = ""
$arr[3,38] = ""
$arr[3,39] = ""
$arr[3,40] = ""

# Row 29
$arr[4,0] = "WGE 84"
$arr[4,1] = "Western Interior Designers & Marine Contractors"
$arr[4,2] = "16-01-2026"
$arr[4,3] = 286962
$arr[4,4] = "Western Interior Designers & Marine Contractors"
$arr[4,5] = 34413429360
$arr[4,6] = "DCR"
$arr[4,7] = "SBIN0003229"
$arr[4,8] = "AAAFW8862C"
$arr[4,9] = "32AAAFW8862C1Z9"
$arr[4,10] = "Lakshmi Priya C B"
$arr[4,11] = "865bfe87-2ceb-4851-ae96-b60d0748fbf0"
$arr[4,12] = 67382981874
$arr[4,13] = "SBIN0012854"
$arr[4,14] = ""
$arr[4,15] = ""
$arr[4,16] = ""
$arr[4,17] = ""
$arr[4,18] = ""
$arr[4,19] = ""
$arr[4,20] = "pending"
$arr[4,21] = 849
$arr[4,22] = ""
$arr[4,23] = "george sir ticket booking expenses ( Ernakulam jn to Madgaon ) credited to lakshmi account RPA_UNIQUE_ID : 86c432fa-489b-49ae-ace8-f85f55a77778"
$arr[4,24] = "ho staff"
$arr[4,25] = 0
$arr[4,26] = "hrm@westernidc.com"
$arr[4,27] = "ESTIMATION NOT MATCHED"
$arr[4,28] = 0
$arr[4,29] = 0
$arr[4,30] = 0
$arr[4,31] = ""
$arr[4,32] = ""
$arr[4,33] = ""
$arr[4,34] = ""
$arr[4,35] = ""
$arr[4,36] = ""
$arr[4,37] = ""
$arr[4,38] = ""
$arr[4,39] = ""
$arr[4,40] = ""

# Row 30
$arr[5,0] = "WGE 73"
$arr[5,1] = "Western Interior Designers & Marine Contractors"
$arr[5,2] = "16-01-2026"
$arr[5,3] = 286962
$arr[5,4] = "Western Interior Designers & Marine Contractors"
$arr[5,5] = 34413429360
$arr[5,6] = "DCR"
$arr[5,7] = "SBIN0003229"
$arr[5,8] = "AAAFW8862C"
$arr[5,9] = "32AAAFW8862C1Z9"
$arr[5,10] = "Nithin"
$arr[5,11] = "0563a8f3-5ed1-4d46-994a-553ef9d83783"
$arr[5,12] = 32555551936
$arr[5,13] = "SBIN0001890"
$arr[5,14] = ""
$arr[5,15] = ""
$arr[5,16] = ""
$arr[5,17] = ""
$arr[5,18] = ""
$arr[5,19] = ""
$arr[5,20] = "pending"
$arr[5,21] = 1000
$arr[5,22] = ""
$arr[5,23] = "Nithin ELATHUR CASUAL WAGE PENDING 31-12=1000`n02-01=1000,03-01=2000,04-01=1000,06-01=1000,07-01=1000`n09-01=2000`nTOTAL =9000 ( ONLY PAID 8000) BALANCE 1000 PENDING RPA_UNIQUE_ID : 5a80203e-c2d2-46e8-9815-af30bbd85df0"
$arr[5,24] = "HPCL ELATHUR"
$arr[5,25] = 0
$arr[5,26] = "hrm@westernidc.com"
$arr[5,27] = "ESTIMATION NOT MATCHED"
$arr[5,28] = 0
$arr[5,29] = 0
$arr[5,30] = 0
$arr[5,31] = ""
$arr[5,32] = ""
$arr[5,33] = ""
$arr[5,34] = ""
$arr[5,35] = ""
$arr[5,36] = ""
$arr[5,37] = ""
$arr[5,38] = ""
$arr[5,39] = ""
$arr[5,40] = ""

# Row 31
$arr[6,0] = "WGE 261"
$arr[6,1] = "Western Interior Designers & Marine Contractors"
$arr[6,2] = "16-01-2026"
$arr[6,3] = 286962
$arr[6,4] = "Western Interior Designers & Marine Contractors"
$arr[6,5] = 34413429360
$arr[6,6] = "NEFT"
$arr[6,7] = "SBIN0003229"
$arr[6,8] = "AAAFW8862C"
$arr[6,9] = "32AAAFW8862C1Z9"
$arr[6,10] = "MS. FATHIMA ( NOUFAL)"
$arr[6,11] = "a56d5450-8b8b-422e-8502-8e313f1948b7"
$arr[6,12] = 337202010014606
$arr[6,13] = "UBIN0533726"
$arr[6,14] = ""
$arr[6,15] = ""
$arr[6,16] = ""
$arr[6,17] = ""
$arr[6,18] = ""
$arr[6,19] = ""
$arr[6,20] = "pending"
$arr[6,21] = 1500
$arr[6,22] = ""
$arr[6,23] = "FOOD EXPENSES( 16/01/2026 TO 31/01/2026) @100 RPA_UNIQUE_ID : bab98b0f-c562-4f03-b3c6-c25c5246ab1a"
$arr[6,24] = "HULL GOA"
$arr[6,25] = 0
$arr[6,26] = "hrm@westernidc.com"
$arr[6,27] = "ESTIMATION NOT MATCHED"
$arr[6,28] = 0
$arr[6,29] = 0
$arr[6,30] = 0
$arr[6,31] = ""
$arr[6,32] = ""
$arr[6,33] = ""
$arr[6,34] = ""
$arr[6,35] = ""
$arr[6,36] = ""
$arr[6,37] = ""
$arr[6,38] = ""
$arr[6,39] = ""
$arr[6,40] = ""

$ws.Range("A25:AO31").Value = $arr
